# Insert a new data row at spreadsheet row 54 (pushes existing rows 54:119 down
# to 55:120, matching the reference diff which shows the whole block re-indexed
# by one row, with a brand-new record appearing as the new row 54 and the
# dimension growing from A1:T119 to A1:T120).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(54).Insert()

$newRow = @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 45033, 16, "Fruta", 100108, "Tropicales y subtropicales", 100108002, "Mango", "Sin especificar", "Primera", 60, 7000, 7000, 7000, "`$/bandeja 4 kilos", "Perú", 1750, 4)

for ($col = 1; $col -le $newRow.Length; $col++) {
    $ws.Cells.Item(54, $col).Value = $newRow[$col - 1]
}
